# Natmi LR-pairs (Myoc-Fzd1): following Dr Hou advice.
# Rebuild the sending/target cluster matrix as a full 3x3 (ECs/FAPs/sCs x ECs/FAPs/sCs)
# instead of the previous 2x3 (FAPs/sCs x ECs/FAPs/sCs), with refreshed NATMI stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Myoc", "Fzd1", "ECs", 2.0, 0.6666666666666666, 0.349434, 1.048302, 0.0009963999680650763, 0.0009963999680650763, 3.0, 1.0, 0.9623583333333334, 2.887075, 0.05805926999654511, 0.05805926999654509, 0.3362807218500001, 3.02652649665, 0.00005785025477043918, 0.00005785025477043916),
  @("ECs", "Myoc", "Fzd1", "FAPs", 2.0, 0.6666666666666666, 0.349434, 1.048302, 0.0009963999680650763, 0.0009963999680650763, 3.0, 1.0, 12.12890133333333, 36.386704, 0.731739034081334, 0.7317390340813339, 4.238250508512, 38.14425457660801, 0.000729104750190611, 0.0007291047501906109),
  @("ECs", "Myoc", "Fzd1", "sCs", 2.0, 0.6666666666666666, 0.349434, 1.048302, 0.0009963999680650763, 0.0009963999680650763, 3.0, 1.0, 3.484187, 10.452561, 0.210201695922121, 0.2102016959221209, 1.217493400158, 10.957440601422, 0.0002094449631040262, 0.0002094449631040261),
  @("FAPs", "Myoc", "Fzd1", "ECs", 3.0, 1.0, 348.977468, 1046.932404, 0.9950981815468188, 0.9950981815468188, 3.0, 1.0, 0.9623583333333334, 2.887075, 0.05805926999654511, 0.05805926999654509, 335.8413744753667, 3022.5723702783, 0.05777467399549781, 0.05777467399549779),
  @("FAPs", "Myoc", "Fzd1", "FAPs", 3.0, 1.0, 348.977468, 1046.932404, 0.9950981815468188, 0.9950981815468188, 3.0, 1.0, 12.12890133333333, 36.386704, 0.731739034081334, 0.7317390340813339, 4232.713276928491, 38094.41949235641, 0.7281521821811612, 0.7281521821811611),
  @("FAPs", "Myoc", "Fzd1", "sCs", 3.0, 1.0, 348.977468, 1046.932404, 0.9950981815468188, 0.9950981815468188, 3.0, 1.0, 3.484187, 10.452561, 0.210201695922121, 0.2102016959221209, 1215.902757298516, 10943.12481568664, 0.2091713253701599, 0.2091713253701599),
  @("sCs", "Myoc", "Fzd1", "ECs", 3.0, 1.0, 1.369616666666667, 4.10885, 0.00390541848511611, 0.00390541848511611, 3.0, 1.0, 0.9623583333333334, 2.887075, 0.05805926999654511, 0.05805926999654509, 1.318062012638889, 11.86255811375, 0.0002267457462768544, 0.0002267457462768543),
  @("sCs", "Myoc", "Fzd1", "FAPs", 3.0, 1.0, 1.369616666666667, 4.10885, 0.00390541848511611, 0.00390541848511611, 3.0, 1.0, 12.12890133333333, 36.386704, 0.731739034081334, 0.7317390340813339, 16.61194541448889, 149.5075087304, 0.002857747149982249, 0.002857747149982249),
  @("sCs", "Myoc", "Fzd1", "sCs", 3.0, 1.0, 1.369616666666667, 4.10885, 0.00390541848511611, 0.00390541848511611, 3.0, 1.0, 3.484187, 10.452561, 0.210201695922121, 0.2102016959221209, 4.772000584983333, 42.94800526485, 0.0008209255888570068, 0.0008209255888570066)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $rows[$i]
  for ($j = 0; $j -lt $r.Length; $j++) {
    $ws.Cells.Item($startRow + $i, $j + 1).Value = $r[$j]
  }
}

